# Fruta / hortaliza, semanal
# A new weekly price record for Mango (Vega Monumental Concepción) is inserted
# as row 153, pushing the existing rows 153-195 down to 154-196.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 153; this shifts rows 153:195 -> 154:196
$ws.Rows.Item(153).Insert()

# Populate the newly inserted row with the latest weekly observation
$ws.Cells.Item(153, 1).Value  = 11
$ws.Cells.Item(153, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(153, 3).Value  = "Bíobío"
$ws.Cells.Item(153, 4).Value  = 45202
$ws.Cells.Item(153, 5).Value  = 8
$ws.Cells.Item(153, 6).Value  = "Fruta"
$ws.Cells.Item(153, 7).Value  = 100108
$ws.Cells.Item(153, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(153, 9).Value  = 100108002
$ws.Cells.Item(153, 10).Value = "Mango"
$ws.Cells.Item(153, 11).Value = "Sin especificar"
$ws.Cells.Item(153, 12).Value = "Primera"
$ws.Cells.Item(153, 13).Value = 200
$ws.Cells.Item(153, 14).Value = 10000
$ws.Cells.Item(153, 15).Value = 11000
$ws.Cells.Item(153, 16).Value = 10400
$ws.Cells.Item(153, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(153, 18).Value = "Brasil"
$ws.Cells.Item(153, 19).Value = 2600
$ws.Cells.Item(153, 20).Value = 4
